$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 4138.364
$ws.Range("I11").Value = 4138.364
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 4138.364
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -3998.364

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 439.725
$ws.Range("I55").Value = 459.34616
$ws.Range("J55").Value = 403.2857
$ws.Range("K55").Value = 459.34616
$ws.Range("L55").Value = 403.2857
$ws.Range("M55").Value = -245.34616
$ws.Range("N55").Value = -831.2857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 216599.8
$ws.Range("I99").Value = 20000
$ws.Range("J99").Value = 265749.75
$ws.Range("K99").Value = 60000
$ws.Range("L99").Value = 797249.25
$ws.Range("M99").Value = -58502
$ws.Range("N99").Value = -800245.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 878.1429000000001
$ws.Range("I132").Value = 819.7659
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 2459.2977
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = 70.70229999999992

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3527.1025
$ws.Range("I137").Value = 2145.6875
$ws.Range("J137").Value = 4488.087
$ws.Range("K137").Value = 6437.0625
$ws.Range("L137").Value = 13464.261
$ws.Range("M137").Value = -3887.0625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4716.6978
$ws.Range("I138").Value = 1326.5714
$ws.Range("J138").Value = 5375.8887
$ws.Range("K138").Value = 3979.7142
$ws.Range("L138").Value = 16127.6661
$ws.Range("M138").Value = 1160.2858
$ws.Range("N138").Value = -26407.6661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1696627.8
$ws.Range("I32").Value = 8929.581
$ws.Range("J32").Value = 12861401
$ws.Range("K32").Value = 8929.581
$ws.Range("L32").Value = 12861401
$ws.Range("M32").Value = -8642.581
$ws.Range("N32").Value = -12861975

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3370.0952
$ws.Range("I61").Value = 2018.8276
$ws.Range("J61").Value = 6384.4614
$ws.Range("K61").Value = 2018.8276
$ws.Range("L61").Value = 6384.4614
$ws.Range("M61").Value = -1806.8276
$ws.Range("N61").Value = -6808.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 13890034
$ws.Range("I102").Value = 849.1
$ws.Range("J102").Value = 83335960
$ws.Range("K102").Value = 849.1
$ws.Range("L102").Value = 83335960
$ws.Range("M102").Value = 772.9
$ws.Range("N102").Value = -83339204

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4558.025
$ws.Range("I132").Value = 3368.3447
$ws.Range("J132").Value = 7694.4546
$ws.Range("K132").Value = 10105.0341
$ws.Range("L132").Value = 23083.3638
$ws.Range("M132").Value = -7575.034100000001
$ws.Range("N132").Value = -28143.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3370.0952
$ws.Range("I136").Value = 2018.8276
$ws.Range("J136").Value = 6384.4614
$ws.Range("K136").Value = 6056.4828
$ws.Range("L136").Value = 19153.3842
$ws.Range("M136").Value = -3506.4828
$ws.Range("N136").Value = -24253.3842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 769.7143
$ws.Range("I16").Value = 769.7143
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 769.7143
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -482.7143
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4567.727
$ws.Range("I31").Value = 3159
$ws.Range("J31").Value = 5976.4546
$ws.Range("K31").Value = 3159
$ws.Range("L31").Value = 5976.4546
$ws.Range("M31").Value = -2864

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4567.727
$ws.Range("I34").Value = 3159
$ws.Range("J34").Value = 5976.4546
$ws.Range("K34").Value = 3159
$ws.Range("L34").Value = 5976.4546
$ws.Range("M34").Value = -2957

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 40000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 40000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 769.7143
$ws.Range("I113").Value = 769.7143
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 769.7143
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1400.2857
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 368219.34
$ws.Range("I122").Value = 682083
$ws.Range("J122").Value = 6069
$ws.Range("K122").Value = 2046249
$ws.Range("L122").Value = 18207
$ws.Range("M122").Value = -2043799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3053
$ws.Range("I132").Value = 2522.3333
$ws.Range("J132").Value = 5918.6
$ws.Range("K132").Value = 7566.999899999999
$ws.Range("L132").Value = 17755.8
$ws.Range("M132").Value = -5036.999899999999
$ws.Range("N132").Value = -22815.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4946
$ws.Range("I5").Value = 566.2778
$ws.Range("J5").Value = 16208.143
$ws.Range("K5").Value = 1698.8334
$ws.Range("L5").Value = 48624.429
$ws.Range("M5").Value = -1586.8334
$ws.Range("N5").Value = -48848.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 252543.3
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 252543.3
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 757629.8999999999
$ws.Range("N37").Value = -757853.8999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1337.3334
$ws.Range("I122").Value = 1446.4
$ws.Range("J122").Value = 1303.25
$ws.Range("K122").Value = 13017.6
$ws.Range("L122").Value = 11729.25
$ws.Range("M122").Value = -10567.6
$ws.Range("N122").Value = -16629.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4521.9414
$ws.Range("I131").Value = 3911.2856
$ws.Range("J131").Value = 4949.4
$ws.Range("K131").Value = 11733.8568
$ws.Range("L131").Value = 14848.2
$ws.Range("M131").Value = -6693.856800000001
$ws.Range("N131").Value = -24928.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 4946
$ws.Range("I135").Value = 566.2778
$ws.Range("J135").Value = 16208.143
$ws.Range("K135").Value = 5096.500199999999
$ws.Range("L135").Value = 145873.287
$ws.Range("M135").Value = -2561.500199999999
$ws.Range("N135").Value = -150943.287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5303.4
$ws.Range("I68").Value = 5174.857
$ws.Range("J68").Value = 5496.2144
$ws.Range("K68").Value = 5174.857
$ws.Range("L68").Value = 5496.2144
$ws.Range("M68").Value = -4425.857
$ws.Range("N68").Value = -6994.2144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 5303.4
$ws.Range("I71").Value = 5174.857
$ws.Range("J71").Value = 5496.2144
$ws.Range("K71").Value = 25874.285
$ws.Range("L71").Value = 27481.072
$ws.Range("M71").Value = -22130.285
$ws.Range("N71").Value = -34969.072

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 41671064
$ws.Range("I81").Value = 2870.5
$ws.Range("J81").Value = 83339256
$ws.Range("K81").Value = 5741
$ws.Range("L81").Value = 166678512
$ws.Range("M81").Value = -4680
$ws.Range("N81").Value = -166680634

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 41671064
$ws.Range("I84").Value = 2870.5
$ws.Range("J84").Value = 83339256
$ws.Range("K84").Value = 28705
$ws.Range("L84").Value = 833392560
$ws.Range("M84").Value = -23401
$ws.Range("N84").Value = -833403168

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 46778.5
$ws.Range("I109").Value = 31342
$ws.Range("J109").Value = 48181.816
$ws.Range("K109").Value = 31342
$ws.Range("L109").Value = 48181.816
$ws.Range("M109").Value = -29955
$ws.Range("N109").Value = -50955.816

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4420.8423
$ws.Range("I122").Value = 3285.4285
$ws.Range("J122").Value = 7600
$ws.Range("K122").Value = 9856.2855
$ws.Range("L122").Value = 22800
$ws.Range("M122").Value = -7406.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3224.4211
$ws.Range("I126").Value = 2476.9333
$ws.Range("J126").Value = 6027.5
$ws.Range("K126").Value = 7430.7999
$ws.Range("L126").Value = 18082.5
$ws.Range("M126").Value = -4960.7999
